# Update Task - 8-1 F/N
#
# Adds a new task row (row 9) to the "sriram" worksheet (2nd tab / the
# tab-selected sheet): S NO=5, DATE=08-Jan-2018, TIME="11:30 to 13:30",
# PROJECT TASK="Completed UiPATH course and obtained Certificate",
# STATUS="completed" styled with the built-in "Good" cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sriram")
$ws.Activate()

# Seed row 9 with row 8's formatting (wrap text on A/C/D, date format on B,
# row height, etc.) so no brand-new near-duplicate style entries get
# minted - only the new "Good" style below is genuinely new.
$ws.Range("A8:E8").Copy()
$ws.Range("A9:E9").PasteSpecial(-4122)

$ws.Range("A9").Value = 5
$ws.Range("B9").Value = Get-Date -Year 2018 -Month 1 -Day 8 -Hour 0 -Minute 0 -Second 0
$ws.Range("C9").Value = "11:30 to 13:30"
$ws.Range("D9").Value = "Completed UiPATH course and obtained Certificate"
$ws.Range("E9").Value = "completed"

# Status cell uses the built-in "Good" (green) cell style, same family as
# row 8's "Neutral" style but without the wrap-text alignment override.
$ws.Range("E9").Style = "Good"

$ws.Rows.Item(9).RowHeight = 30

$ws.Range("E10").Select()
